$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "_200001_ReadBookA"
$ws.Range("B4").Value = "_200002_HitSuspiciousTrigger"
$ws.Range("C6").Select()
